# Refresh the cryptocurrency price/volume snapshot (GitHub Actions scheduled update).
# Column D ("Price") and column E ("Volume(1h)") are plain text cells (not numbers),
# so numeric-looking prices are written with NumberFormat "@" (Text) first, then the
# cell style is reset to Normal so no stray number-format style lingers on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D2='69.187.00', E2='  +2.39%  '
$ws.Range('D2').Value = '69.187.00'
$ws.Range('E2').Value = '  +2.39%  '

# Row 3: D3='3.815.09', E3='  +0.89%  '
$ws.Range('D3').Value = '3.815.09'
$ws.Range('E3').Value = '  +0.89%  '

# Row 4: D4='0.993', E4='  -0.88%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.993'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.88%  '

# Row 5: D5='629.89', E5='  +5.31%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '629.89'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.31%  '

# Row 6: D6='165.07', E6='  +0.35%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '165.07'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.35%  '

# Row 7: D7='3.813.06', E7='  +0.92%  '
$ws.Range('D7').Value = '3.813.06'
$ws.Range('E7').Value = '  +0.92%  '

# Row 8: D8='0.999', E8='  -0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.06%  '

# Row 9: D9='0.520', E9='  +1.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.520'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.10%  '

# Row 10: E10='  +2.84%  '
$ws.Range('E10').Value = '  +2.84%  '

# Row 11: D11='0.454', E11='  +1.07%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.454'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.07%  '

# Row 12: D12='6.59', E12='  +3.06%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.59'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.06%  '

# Row 13: E13='  +1.06%  '
$ws.Range('E13').Value = '  +1.06%  '

# Row 14: D14='35.99', E14='  +1.43%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.99'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.43%  '

# Row 15: D15='4.449.16', E15='  +0.79%  '
$ws.Range('D15').Value = '4.449.16'
$ws.Range('E15').Value = '  +0.79%  '

# Row 16: D16='3.805.70', E16='  +0.71%  '
$ws.Range('D16').Value = '3.805.70'
$ws.Range('E16').Value = '  +0.71%  '

# Row 17: D17='69.104.02', E17='  +2.24%  '
$ws.Range('D17').Value = '69.104.02'
$ws.Range('E17').Value = '  +2.24%  '

# Row 18: D18='17.95', E18='  -1.95%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.95'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.95%  '

# Row 19: D19='7.13', E19='  +1.31%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.13'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.31%  '

# Row 20: E20='  -0.13%  '
$ws.Range('E20').Value = '  -0.13%  '

# Row 21: D21='466.18', E21='  +1.29%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '466.18'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.29%  '

# Row 22: D22='9.67', E22='  -0.28%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.67'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.28%  '

# Row 23: E23='  +2.02%  '
$ws.Range('E23').Value = '  +2.02%  '

# Row 24: E24='  +3.72%  '
$ws.Range('E24').Value = '  +3.72%  '

# Row 25: D25='83.66', E25='  +1.40%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.66'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.40%  '

# Row 26: E26='  +3.62%  '
$ws.Range('E26').Value = '  +3.62%  '

# Row 27: D27='11.98', E27='  +0.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.98'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.09%  '

# Row 28: E28='  +1.64%  '
$ws.Range('E28').Value = '  +1.64%  '

# Row 29: E29='  +0.03%  '
$ws.Range('E29').Value = '  +0.03%  '

# Row 30: D30='3.959.61', E30='  +0.81%  '
$ws.Range('D30').Value = '3.959.61'
$ws.Range('E30').Value = '  +0.81%  '

# Row 31: D31='2.70', E31='  +3.83%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.70'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.83%  '

# Row 33: D33='7.28', E33='  -1.92%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.28'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.92%  '

# Row 34: D34='29.19', E34='  +0.55%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '29.19'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.55%  '

# Row 35: D35='0.999', E35='  -0.13%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.13%  '

# Row 36: E36='  +1.41%  '
$ws.Range('E36').Value = '  +1.41%  '

# Row 37: D37='0.103', E37='  +3.80%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.103'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.80%  '

# Row 38: D38='0.150', E38='  +8.21%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.150'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +8.21%  '

# Row 39: D39='3.48', E39='  +7.04%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.48'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +7.04%  '

# Row 40: E40='  +2.72%  '
$ws.Range('E40').Value = '  +2.72%  '

# Row 41: D41='0.976', E41='  -0.84%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.976'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.84%  '

# Row 42: E42='  +0.06%  '
$ws.Range('E42').Value = '  +0.06%  '

# Row 44: D44='157.45', E44='  +3.70%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '157.45'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.70%  '

# Row 45: E45='  +1.43%  '
$ws.Range('E45').Value = '  +1.43%  '

# Row 46: E46='  +5.19%  '
$ws.Range('E46').Value = '  +5.19%  '

# Row 47: D47='43.31', E47='  -0.14%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '43.31'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.14%  '

# Row 48: D48='46.90', E48='  -1.15%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '46.90'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.15%  '

# Row 49: E49='  +3.54%  '
$ws.Range('E49').Value = '  +3.54%  '

# Row 50: D50='8.44', E50='  +1.51%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.44'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.51%  '

# Row 51: D51='0.000280', E51='  +13.70%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.000280'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +13.70%  '
